# "download articles with pandoc title blocks"
#
# Pandoc's docx writer now emits a proper title block (Title / Authors
# paragraph styles, one run per word/space/punctuation "token" instead of
# a single run per paragraph) rather than a Heading1 title followed by a
# bold "By <author>" paragraph.
#
# This script rebuilds the first two paragraphs of the document to match
# that shape:
#   1) "On Pilgrimage - March-April 1970" (Heading1, one run) becomes a
#      Title-styled paragraph, split word-by-word/space-by-space into
#      separate runs.
#   2) "By Dorothy Day" (bold, one run) becomes an Authors-styled
#      paragraph containing just "Dorothy Day", likewise split into
#      separate runs, with the "By " prefix and bold formatting dropped.

$d = $word.ActiveDocument

function New-WordRunsXml([string[]]$tokens) {
    $sb = New-Object System.Text.StringBuilder
    foreach ($tok in $tokens) {
        $escaped = $tok.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
        [void]$sb.Append("<w:r><w:t xml:space=`"preserve`">$escaped</w:t></w:r>")
    }
    return $sb.ToString()
}

function New-PandocParagraphXml([string]$style, [string[]]$tokens) {
    $runs = New-WordRunsXml $tokens
    return "<w:p><w:pPr><w:pStyle w:val=`"$style`"/></w:pPr>$runs</w:p>"
}

$titleTokens = @("On", " ", "Pilgrimage", " ", "-", " ", "March", "-", "April", " ", "1970")
$authorTokens = @("Dorothy", " ", "Day")

$titleParaXml = New-PandocParagraphXml "Title" $titleTokens
$authorsParaXml = New-PandocParagraphXml "Authors" $authorTokens

$bodyXml = "$titleParaXml$authorsParaXml"

$packageXml = "<?xml version=`"1.0`" encoding=`"UTF-8`" standalone=`"yes`"?>" +
    "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" +
    "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" +
    "<pkg:xmlData>" +
    "<w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">" +
    "<w:body>$bodyXml</w:body></w:document>" +
    "</pkg:xmlData></pkg:part></pkg:package>"

# Paragraph 1 is the "On Pilgrimage - March-April 1970" Heading1 title;
# paragraph 2 is the bold "By Dorothy Day" byline. Replace both (including
# their paragraph marks) in one shot so the new Title/Authors paragraphs
# land exactly where the old ones were.
$titlePara = $d.Paragraphs(1)
$authorPara = $d.Paragraphs(2)
$target = $d.Range($titlePara.Range.Start, $authorPara.Range.End)
$target.InsertXML($packageXml)

# Best-effort: also drop the old "on-pilgrimage---march-april-1970"
# bookmark that used to wrap the title paragraph (pandoc's new output no
# longer emits it). Some runtimes don't surface malformed/legacy
# body-level bookmarks through the Bookmarks collection, so this is
# wrapped defensively and is not relied upon.
try {
    for ($i = $d.Bookmarks.Count; $i -ge 1; $i--) {
        $bm = $d.Bookmarks.Item($i)
        if ($bm.Name -eq "on-pilgrimage---march-april-1970") {
            $bm.Delete()
        }
    }
} catch {
}

Write-Output $d.Paragraphs(1).Range.Text
Write-Output $d.Paragraphs(1).Style
Write-Output $d.Paragraphs(2).Range.Text
Write-Output $d.Paragraphs(2).Style
